$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The project's raw data collection was completed and the workbook was
# trimmed down to the 11 Southeast Asian countries actually in scope
# (commit: "Complete mobile data retrieved. Limited data set to SE Asia
# 11 countries."). Remove the rows for the countries that fell outside
# that SE Asia scope: Bangladesh, Iraq, Maldives, Mongolia, Pakistan and
# Vanuatu. Delete bottom-up so earlier row numbers stay valid.
$ws.Rows("17:17").Delete() | Out-Null   # Vanuatu
$ws.Rows("12:12").Delete() | Out-Null   # Pakistan
$ws.Rows("9:10").Delete()  | Out-Null   # Maldives, Mongolia
$ws.Rows("6:6").Delete()   | Out-Null   # Iraq
$ws.Rows("2:2").Delete()   | Out-Null   # Bangladesh

# Re-assert the formula for the new row 2 (now Brunei Darussalam) so it
# round-trips cleanly now that it sits next to the shared-formula block.
$ws.Range("H2").Formula = "=CEILING.MATH(G2)"

# Update the active selection left behind on the sheet.
$ws.Range("C7").Select() | Out-Null
